$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two rows that were dropped from the dataset entirely
# (original row 26 = "RM 232", original row 28 = "SC 92").
# Deleting row 26 first shifts "SC 92" up to row 27.
$ws.Rows(26).Delete()
$ws.Rows(27).Delete()

# After the two deletions the remaining rows have shifted up by one or
# two, and a handful of cells swap between a value and "missing" (blank).

# F5 (RM 14): value -> missing
$ws.Range("F5").ClearContents()

# F11 (RM 58): missing -> value
$ws.Range("F11").Value = 17.65

# E19 (RM 125): missing -> value ; F19: value -> missing
$ws.Range("E19").Value = -6.5
$ws.Range("F19").ClearContents()

# E21 (RM 135): value -> missing
$ws.Range("E21").ClearContents()

# E23 (RM 140): missing -> value
$ws.Range("E23").Value = -7

# F25 (RM 145): missing -> value
$ws.Range("F25").Value = 16.6

# E27 (SC 101, was row 29 before deletions): value -> missing
$ws.Range("E27").ClearContents()

# F29 (SC 119, was row 31 before deletions): value -> missing
$ws.Range("F29").ClearContents()

# E33 (SC 232, was row 35 before deletions): missing -> value
$ws.Range("E33").Value = -10.7
